# Apply the diff: add a new ListParagraph-styled bullet after the
# "So the main constraints..." paragraph, describing the initial
# solution, and move the trailing _GoBack bookmark onto the new
# paragraph (matching Word's behaviour of tracking the last edit spot).

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14"'

function New-FlatOpcXml($bodyInner) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $wNs + '><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData>' +
        '</pkg:part></pkg:package>'
}

# Locate the "So the main constraints..." paragraph (third paragraph).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "So the main constraints to this problem*") {
        $targetPara = $cand
        break
    }
}

$apos = [char]0x2019

# Rewrite that paragraph's content without the _GoBack bookmark (the
# bookmark will be re-created on the new paragraph below, mirroring
# where Word leaves it after the last edit).
$existingText = "So the main constraints to this problem are that the cat and parrot can" + $apos + "t be left alone. The parrot and bag of seed can" + $apos + "t be left alone. So the goal will be to make sure each item gets across without the constraints being applied."

$existingParaBody = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>' + $existingText + '</w:t></w:r>'

$targetRange = $targetPara.Range
$targetRange.InsertXML((New-FlatOpcXml ('<w:p>' + $existingParaBody + '</w:p>'))) | Out-Null

# Re-fetch the paragraph and insert a new paragraph right after it.
$targetPara = $d.Paragraphs.Item(3)
$endRange = $targetPara.Range
$endRange.Collapse(0) | Out-Null
$endRange.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item(4)

$newParaBody = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>The initial solution would be to take the parrot across first. Return empty handed and pick up the cat in order to take it to the other side leaving the bag of seed alone. When getting to the other side, the man drop</w:t></w:r>' +
    '<w:r><w:t>s</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> the cat off while picking up the parrot to return it to the other side so </w:t></w:r>' +
    '<w:r><w:t>the cat does not eat it. When back on the initial side the man drops off the parrot while picking up the bag of seed to take to the other side. Now the man can drop off the seed return to pick up the parrot and will have successfully transported all items across the riverbank.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$newRange = $newPara.Range
$newRange.InsertXML((New-FlatOpcXml ('<w:p>' + $newParaBody + '</w:p>'))) | Out-Null

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
